$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 24-31 need the same formatting (bold + border + centered) that
# column A carries on every existing data row (style index 1 / row 3 here).
$ws.Range("A3").Copy()
$ws.Range("A24:A31").PasteSpecial(-4122)  # xlPasteFormats

# Update process numbers (rows 3-23 changed values; rows 24-31 are brand new)

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "5002825-34.2018.8.21.0002"
$ws.Cells.Item(3, 3).Value = "9000482-94.2018.8.21.0002"
$ws.Cells.Item(3, 4).Value = "CIV.12883.01"
$ws.Cells.Item(3, 5).Value = "originario_principal"

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "5010408-92.2022.8.21.0014"
$ws.Cells.Item(4, 3).Value = "5002665-02.2020.8.21.0014"
$ws.Cells.Item(4, 4).Value = "CIV.36852.01"
$ws.Cells.Item(4, 5).Value = "originario_principal"

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "5034684-18.2011.8.21.0001"
$ws.Cells.Item(5, 3).Value = "0111143-49.2011.8.21.0001"
$ws.Cells.Item(5, 4).Value = "CIV.16263.01"
$ws.Cells.Item(5, 5).Value = "originario_principal"

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "5029337-62.2015.8.21.0001"
$ws.Cells.Item(6, 3).Value = "0090773-10.2015.8.21.0001"
$ws.Cells.Item(6, 4).Value = "CIV.00672.01"
$ws.Cells.Item(6, 5).Value = "originario_principal"

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "5039775-79.2017.8.21.0001"
$ws.Cells.Item(7, 3).Value = "0039384-15.2017.8.21.0001"
$ws.Cells.Item(7, 4).Value = "CIV.17971.01"
$ws.Cells.Item(7, 5).Value = "originario_principal"

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "5029499-57.2015.8.21.0001"
$ws.Cells.Item(8, 3).Value = "0136341-49.2015.8.21.0001"
$ws.Cells.Item(8, 4).Value = "CIV.18618.01"
$ws.Cells.Item(8, 5).Value = "originario_principal"

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "5028728-21.2011.8.21.0001"
$ws.Cells.Item(9, 3).Value = "0344069-02.2011.8.21.0001"
$ws.Cells.Item(9, 4).Value = "CIV.18456.01"
$ws.Cells.Item(9, 5).Value = "originario_principal"

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "5028728-21.2011.8.21.0001"
$ws.Cells.Item(10, 3).Value = "0344069-02.2011.8.21.0001"
$ws.Cells.Item(10, 4).Value = "CIV.18456.01"
$ws.Cells.Item(10, 5).Value = "originario_principal"

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "5001806-14.2020.8.21.0134"
$ws.Cells.Item(11, 3).Value = "9000364-42.2020.8.21.0134"
$ws.Cells.Item(11, 4).Value = "CIV.35818.01"
$ws.Cells.Item(11, 5).Value = "originario_principal"

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "5032111-07.2011.8.21.0001"
$ws.Cells.Item(12, 3).Value = "0111095-90.2011.8.21.0001"
$ws.Cells.Item(12, 4).Value = "CIV.18452.01"
$ws.Cells.Item(12, 5).Value = "originario_principal"

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "5033921-46.2013.8.21.0001"
$ws.Cells.Item(13, 3).Value = "0037413-34.2013.8.21.0001"
$ws.Cells.Item(13, 4).Value = "CIV.18689.01"
$ws.Cells.Item(13, 5).Value = "originario_principal"

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "5033938-82.2013.8.21.0001"
$ws.Cells.Item(14, 3).Value = "0084903-52.2013.8.21.0001"
$ws.Cells.Item(14, 4).Value = "CIV.18585.01"
$ws.Cells.Item(14, 5).Value = "originario_principal"

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "5033938-82.2013.8.21.0001"
$ws.Cells.Item(15, 3).Value = "0084903-52.2013.8.21.0001"
$ws.Cells.Item(15, 4).Value = "CIV.18585.01"
$ws.Cells.Item(15, 5).Value = "originario_principal"

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "5034971-78.2011.8.21.0001"
$ws.Cells.Item(16, 3).Value = "0108702-95.2011.8.21.0001"
$ws.Cells.Item(16, 4).Value = "CIV.08105.01"
$ws.Cells.Item(16, 5).Value = "originario_principal"

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "5034971-78.2011.8.21.0001"
$ws.Cells.Item(17, 3).Value = "0108702-95.2011.8.21.0001"
$ws.Cells.Item(17, 4).Value = "CIV.08105.01"
$ws.Cells.Item(17, 5).Value = "originario_principal"

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "5039859-80.2017.8.21.0001"
$ws.Cells.Item(18, 3).Value = "0129049-42.2017.8.21.0001"
$ws.Cells.Item(18, 4).Value = "CIV.18176.01"
$ws.Cells.Item(18, 5).Value = "originario_principal"

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "5015478-52.2010.8.21.0001"
$ws.Cells.Item(19, 3).Value = "2626731-39.2010.8.21.0001"
$ws.Cells.Item(19, 4).Value = "CIV.37759.01"
$ws.Cells.Item(19, 5).Value = "originario_principal"

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "5015600-65.2010.8.21.0001"
$ws.Cells.Item(20, 3).Value = "0155171-39.2010.8.21.0001"
$ws.Cells.Item(20, 4).Value = "CIV.08455.01"
$ws.Cells.Item(20, 5).Value = "originario_principal"

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "5034891-17.2011.8.21.0001"
$ws.Cells.Item(21, 3).Value = "0243442-87.2011.8.21.0001"
$ws.Cells.Item(21, 4).Value = "CIV.11451.01"
$ws.Cells.Item(21, 5).Value = "originario_principal"

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = "5034891-17.2011.8.21.0001"
$ws.Cells.Item(22, 3).Value = "0243442-87.2011.8.21.0001"
$ws.Cells.Item(22, 4).Value = "CIV.11451.01"
$ws.Cells.Item(22, 5).Value = "originario_principal"

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "5001001-39.2018.8.21.0067"
$ws.Cells.Item(23, 3).Value = "9000812-90.2018.8.21.0067"
$ws.Cells.Item(23, 4).Value = "CIV.35926.01"
$ws.Cells.Item(23, 5).Value = "originario_principal"

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "5002047-92.2020.8.21.0067"
$ws.Cells.Item(24, 3).Value = "9000436-36.2020.8.21.0067"
$ws.Cells.Item(24, 4).Value = "CIV.36474.01"
$ws.Cells.Item(24, 5).Value = "originario_principal"

$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "5002099-25.2019.8.21.0067"
$ws.Cells.Item(25, 3).Value = "9000760-60.2019.8.21.0067"
$ws.Cells.Item(25, 4).Value = "CIV.35962.01"
$ws.Cells.Item(25, 5).Value = "originario_principal"

$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "5001002-24.2018.8.21.0067"
$ws.Cells.Item(26, 3).Value = "9000850-05.2018.8.21.0067"
$ws.Cells.Item(26, 4).Value = "CIV.35910.01"
$ws.Cells.Item(26, 5).Value = "originario_principal"

$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "5002100-10.2019.8.21.0067"
$ws.Cells.Item(27, 3).Value = "9000786-58.2019.8.21.0067"
$ws.Cells.Item(27, 4).Value = "CIV.35964.01"
$ws.Cells.Item(27, 5).Value = "originario_principal"

$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "5002101-92.2019.8.21.0067"
$ws.Cells.Item(28, 3).Value = "9000444-47.2019.8.21.0067"
$ws.Cells.Item(28, 4).Value = "CIV.06132.01"
$ws.Cells.Item(28, 5).Value = "originario_principal"

$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "5002102-77.2019.8.21.0067"
$ws.Cells.Item(29, 3).Value = "9000756-23.2019.8.21.0067"
$ws.Cells.Item(29, 4).Value = "CIV.35898.01"
$ws.Cells.Item(29, 5).Value = "originario_principal"

$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "5002103-62.2019.8.21.0067"
$ws.Cells.Item(30, 3).Value = "9000754-53.2019.8.21.0067"
$ws.Cells.Item(30, 4).Value = "CIV.35897.01"
$ws.Cells.Item(30, 5).Value = "originario_principal"

$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "5002104-47.2019.8.21.0067"
$ws.Cells.Item(31, 3).Value = "9000700-87.2019.8.21.0067"
$ws.Cells.Item(31, 4).Value = "CIV.35974.01"
$ws.Cells.Item(31, 5).Value = "originario_principal"
